# Update Leve profit figures (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) to reflect refreshed Universalis market data.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 6639
$ws.Range("I74").Value = 3854.875
$ws.Range("J74").Value = 7876.3887
$ws.Range("K74").Value = 3854.875
$ws.Range("L74").Value = 7876.3887
$ws.Range("M74").Value = -2918.875
$ws.Range("N74").Value = -9748.3887

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 6639
$ws.Range("I77").Value = 3854.875
$ws.Range("J77").Value = 7876.3887
$ws.Range("K77").Value = 19274.375
$ws.Range("L77").Value = 39381.9435
$ws.Range("M77").Value = -14594.375
$ws.Range("N77").Value = -48741.9435

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 2453.2727
$ws.Range("I86").Value = 2291.8572
$ws.Range("K86").Value = 2291.8572
$ws.Range("M86").Value = -1168.8572

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 2453.2727
$ws.Range("I89").Value = 2291.8572
$ws.Range("K89").Value = 11459.286
$ws.Range("M89").Value = -5843.286

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2399.2334
$ws.Range("J138").Value = 4737.091
$ws.Range("L138").Value = 14211.273
$ws.Range("N138").Value = -24491.273

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 5601.2075
$ws.Range("I32").Value = 3451.15
$ws.Range("K32").Value = 3451.15
$ws.Range("M32").Value = -3164.15

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2255.1428
$ws.Range("I61").Value = 2094.5454
$ws.Range("J61").Value = 2526.923
$ws.Range("K61").Value = 2094.5454
$ws.Range("L61").Value = 2526.923
$ws.Range("M61").Value = -1882.5454
$ws.Range("N61").Value = -2950.923

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 580947
$ws.Range("I122").Value = 2019.5769
$ws.Range("J122").Value = 2086158.2
$ws.Range("K122").Value = 6058.7307
$ws.Range("L122").Value = 6258474.6
$ws.Range("M122").Value = -3608.7307
$ws.Range("N122").Value = -6263374.6

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2271.6667
$ws.Range("I132").Value = 2028.909
$ws.Range("J132").Value = 2653.1428
$ws.Range("K132").Value = 6086.727000000001
$ws.Range("L132").Value = 7959.428400000001
$ws.Range("M132").Value = -3556.727000000001
$ws.Range("N132").Value = -13019.4284

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2255.1428
$ws.Range("I136").Value = 2094.5454
$ws.Range("J136").Value = 2526.923
$ws.Range("K136").Value = 6283.6362
$ws.Range("L136").Value = 7580.768999999999
$ws.Range("M136").Value = -3733.6362
$ws.Range("N136").Value = -12680.769

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 4468291.5
$ws.Range("I99").Value = 6214866.5
$ws.Range("J99").Value = 4821.1113
$ws.Range("K99").Value = 6214866.5
$ws.Range("L99").Value = 4821.1113
$ws.Range("M99").Value = -6213368.5
$ws.Range("N99").Value = -7817.1113

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3482.1538
$ws.Range("I58").Value = 2853.8572
$ws.Range("J58").Value = 4215.1665
$ws.Range("K58").Value = 2853.8572
$ws.Range("L58").Value = 4215.1665
$ws.Range("M58").Value = -2650.8572
$ws.Range("N58").Value = -4621.1665

# Row 99: O Pine
$ws.Range("H99").Value = 4364
$ws.Range("I99").Value = 4128
$ws.Range("J99").Value = 4600
$ws.Range("K99").Value = 4128
$ws.Range("L99").Value = 4600
$ws.Range("M99").Value = -2630
$ws.Range("N99").Value = -7596

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 3308.111
$ws.Range("I122").Value = 2803.8333
$ws.Range("J122").Value = 4316.6665
$ws.Range("K122").Value = 8411.499899999999
$ws.Range("L122").Value = 12949.9995
$ws.Range("M122").Value = -5961.499899999999
$ws.Range("N122").Value = -17849.9995

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4364
$ws.Range("I126").Value = 4128
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 12384
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -9914
$ws.Range("N126").Value = -18740

# Row 136: Turali Quality
$ws.Range("H136").Value = 3482.1538
$ws.Range("I136").Value = 2853.8572
$ws.Range("J136").Value = 4215.1665
$ws.Range("K136").Value = 8561.5716
$ws.Range("L136").Value = 12645.4995
$ws.Range("M136").Value = -6011.571599999999
$ws.Range("N136").Value = -17745.4995

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 250
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 133: Friends Are Food
$ws.Range("H133").Value = 1230
$ws.Range("I133").Value = 1230
$ws.Range("K133").Value = 3690
$ws.Range("M133").Value = 1370

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 10533152
$ws.Range("I70").Value = 15391684
$ws.Range("J70").Value = 6333
$ws.Range("K70").Value = 15391684
$ws.Range("L70").Value = 6333
$ws.Range("M70").Value = -15391414
$ws.Range("N70").Value = -6873

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 10533152
$ws.Range("I73").Value = 15391684
$ws.Range("J73").Value = 6333
$ws.Range("K73").Value = 15391684
$ws.Range("L73").Value = 6333
$ws.Range("M73").Value = -15390748
$ws.Range("N73").Value = -8205

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2221249.2
$ws.Range("J80").Value = 8366.667
$ws.Range("L80").Value = 8366.667
$ws.Range("N80").Value = -10362.667

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2221249.2
$ws.Range("J83").Value = 8366.667
$ws.Range("L83").Value = 41833.335
$ws.Range("N83").Value = -51817.335

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 33336694
$ws.Range("I113").Value = 83334340
$ws.Range("J113").Value = 4933.3335
$ws.Range("K113").Value = 83334340
$ws.Range("L113").Value = 4933.3335
$ws.Range("M113").Value = -83332170
$ws.Range("N113").Value = -9273.3335

# Row 132: On Board for Lar
$ws.Range("H132").Value = 4058
$ws.Range("I132").Value = 3231.875
$ws.Range("J132").Value = 5526.6665
$ws.Range("K132").Value = 9695.625
$ws.Range("L132").Value = 16579.9995
$ws.Range("M132").Value = -7165.625
$ws.Range("N132").Value = -21639.9995

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 5000
$ws.Range("M68").Value = -4251

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = -21256

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 23811930
$ws.Range("I93").Value = 25643502
$ws.Range("K93").Value = 25643502
$ws.Range("M93").Value = -25642254

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 9726.305
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 9850.228
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 9850.228
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -11098.228

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 9726.305
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 9850.228
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 49251.14
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -55491.14

# Row 69: Fashion Patrol
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31498

# Row 72: Dress Code Violation (L)
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -97488
